$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 98, shifting existing rows 98-131 down to 99-132.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly entry.
$ws.Cells.Item(98, 1).Value = 10
$ws.Cells.Item(98, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(98, 3).Value = "La Araucanía"
$ws.Cells.Item(98, 4).Value = 45135
$ws.Cells.Item(98, 5).Value = 9
$ws.Cells.Item(98, 6).Value = 100112010
$ws.Cells.Item(98, 7).Value = "Achicoria"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 75
$ws.Cells.Item(98, 11).Value = 10000
$ws.Cells.Item(98, 12).Value = 10000
$ws.Cells.Item(98, 13).Value = 10000
$ws.Cells.Item(98, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(98, 15).Value = "Región Metropolitana"
$ws.Cells.Item(98, 16).Value = 556
$ws.Cells.Item(98, 17).Value = 18
$ws.Cells.Item(98, 18).Value = "Hortaliza"
